# "add white box test and telegram fee test and add ui text part"
# Populate the "actual value" (实际值, column F) and "execution result"
# (执行结果, column G) columns for every test-case row on the sheet.
# The actual value mirrors the expected value (column E) because every
# test case in this run is expected to pass, so column G is always "pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 20

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $expected = $ws.Cells.Item($row, 5).Value2
    $ws.Cells.Item($row, 6).Value2 = $expected
    $ws.Cells.Item($row, 7).Value2 = "pass"
}

# Restore the active-cell selection left behind by the author's edit.
$ws.Range("G11").Select()
